# CDUN - (Administrar Personal) Modificar Hostel Worker.xlsx
# Update narrative use case: extend precondition text, reorder the normal
# course steps (1-6), renumber the alternate course labels, and reset the
# active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Precondición: add the "must be Administrator" clause
$ws.Range("E6").Value = "Previa existencia del Hostel Worker en el sistema. El usuario debe ser Administrador."

# Curso normal (Básico) - reorder / rewrite the six narrative steps
$ws.Range("C17").Value = "El caso de uso comienza cuando el administrador selecciona Modificar Hostel Worker"
$ws.Range("G17").Value = "El sistema muestra los datos modificables del sistema."
$ws.Range("C18").Value = "El administrador selecciona el Hostel Worker a modificar."
$ws.Range("G18").Value = "El administrador modifica los datos del Hostel Worker."
$ws.Range("C19").Value = "El sistema comprueba la validez de los datos y los guarda en el sistema"
$ws.Range("F19").Value = 6
$ws.Range("G19").Value = "Fin del caso de uso"

# Cursos alternos - renumber labels
$ws.Range("B22").Value = "2a"
$ws.Range("B23").Value = "5a"

# Restore the view selection state
$ws.Range("K8").Select()
